# Bring in DeleteAsmtTest case.
# Replace the old sample QA-assignment-create data with the new
# DeleteAsmtTest dataset (wider set of columns incl. questions /
# post-submission-instructions, plus a long free-text repro note).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "name_input"
$ws.Range("C1").Value = "due_date"
$ws.Range("D1").Value = "grade_type"
$ws.Range("E1").Value = "instructions"
$ws.Range("F1").Value = "post_sub_instructions"
$ws.Range("G1").Value = "qestion1"
$ws.Range("H1").Value = "question2"
$ws.Range("I1").Value = "question3"

# ---- Row 2 ------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Ren_QA_0001"
# C2 keeps its existing due_date (4/23/2018) -- unchanged by this edit.
$ws.Range("D2").Value = "Pass/Fail"
$ws.Range("E2").Value = 'ins_symbol~!@#$%^&*()_+=-0987654321`{}:"|<>?][' + "'" + ';,./'
$ws.Range("F2").Value = "This message is post submission instructions text."
$ws.Range("G2").Value = "what's tcp/ip protocol?"
$ws.Range("H2").Value = "reason about implemented error"

# ---- Row 3 ------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Ren_QA_0002"
# C3 keeps its existing due_date (4/29/2018) -- unchanged by this edit.
$ws.Range("D3").Value = "Rubric"
$ws.Range("E3").Value = "This message is instructions text."
$ws.Range("F3").Value = "PS_in Test long string less than 400 charactorsSteps to reproduce:1. Login bigben (https:bigben-moodle.youseeu.com) as educator-1.2. Select Course - bigbengenerallink.3. Create any type of project assignment, such as individual project and save it.4. Edit this assignment.5. Click the revord icon in Instructions.6. Click UPLOAD VIDEO button.7. Click SELECT FILE button to select file from local host."

# Row 3 needs to be taller to show the long instructions text, and F3 (plus
# the placeholder cell down at F12) get word-wrap so the text is readable.
$ws.Range("F3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 90
$ws.Range("F12").WrapText = $true

# ---- Column widths for the new / widened columns ----------------------
$ws.Columns.Item(5).ColumnWidth = 50.14
$ws.Columns.Item(6).ColumnWidth = 70.71
$ws.Columns.Item(7).ColumnWidth = 19.86
$ws.Columns.Item(8).ColumnWidth = 30

# ---- Selection, matching the saved workbook's last active cell --------
$ws.Range("F3").Select() | Out-Null
